$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.313.38'
$ws.Range("E2").Value = '  +1.15%  '

$ws.Range("D3").Value = '''1.810.60'
$ws.Range("E3").Value = '  +3.39%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''338.08'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").Value = '''0.9990'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").Value = '''0.4657'
$ws.Range("E7").Value = '  +20.70%  '

$ws.Range("D8").Value = '''0.3815'
$ws.Range("E8").Value = '  +12.27%  '

$ws.Range("D9").Value = '''45.56'
$ws.Range("E9").Value = '  -0.50%  '

$ws.Range("D10").Value = '''1.161'
$ws.Range("E10").Value = '  +3.72%  '

$ws.Range("D11").Value = '''0.07659'
$ws.Range("E11").Value = '  +6.04%  '

$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("E13").Value = '  -0.06%  '

$ws.Range("D14").Value = '''6.356'
$ws.Range("E14").Value = '  +3.10%  '

$ws.Range("D15").Value = '''7.462'
$ws.Range("E15").Value = '  +4.84%  '

$ws.Range("D16").Value = '''1.805.90'
$ws.Range("E16").Value = '  +3.04%  '

$ws.Range("E17").Value = '  +3.43%  '

$ws.Range("D18").Value = '''0.06717'
$ws.Range("E18").Value = '  +1.67%  '

$ws.Range("D19").Value = '''81.98'
$ws.Range("E19").Value = '  +3.44%  '

$ws.Range("D20").Value = '''0.9991'
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = '''17.53'
$ws.Range("E21").Value = '  +4.61%  '

$ws.Range("D22").Value = '''6.436'
$ws.Range("E22").Value = '  +4.12%  '

$ws.Range("D23").Value = '''28.305.99'

$ws.Range("E24").Value = '  +2.25%  '

$ws.Range("D25").Value = '''2.416'
$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("D26").Value = '''20.83'
$ws.Range("E26").Value = '  +4.83%  '

$ws.Range("D27").Value = '''154.02'
$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("D28").Value = '''2.386'
$ws.Range("E28").Value = '  +3.81%  '

$ws.Range("D29").Value = '''2.014.83'
$ws.Range("E29").Value = '  +3.24%  '

$ws.Range("E30").Value = '  +1.99%  '

$ws.Range("D31").Value = '''1.263'
$ws.Range("E31").Value = '  -0.78%  '

$ws.Range("D32").Value = '''4.037'
$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").Value = '''0.09590'
$ws.Range("E33").Value = '  +8.61%  '

$ws.Range("D34").Value = '''5.879'
$ws.Range("E34").Value = '  +0.72%  '

$ws.Range("D35").Value = '''0.2281'
$ws.Range("E35").Value = '  +8.72%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.06399'
$ws.Range("E36").Value = '  +4.12%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '''12.15'
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").Value = '''0.02362'

$ws.Range("D39").Value = '''5.294'
$ws.Range("E39").Value = '  +3.03%  '

$ws.Range("D40").Value = '''0.6667'
$ws.Range("E40").Value = '  +1.65%  '

$ws.Range("D41").Value = '''1.242'
$ws.Range("E41").Value = '  +2.89%  '

$ws.Range("D42").Value = '''1.494'
$ws.Range("E42").Value = '  -3.01%  '

$ws.Range("D43").Value = '''8.338'
$ws.Range("E43").Value = '  +3.91%  '

$ws.Range("D44").Value = '''14.27'
$ws.Range("E44").Value = '  +4.30%  '

$ws.Range("D45").Value = '''0.9988'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '''0.6164'
$ws.Range("E46").Value = '  +2.10%  '

$ws.Range("E47").Value = '  +0.48%  '

$ws.Range("D48").Value = '''131.10'
$ws.Range("E48").Value = '  +3.47%  '

$ws.Range("D49").Value = '''2.047'
$ws.Range("E49").Value = '  +2.11%  '

$ws.Range("E50").Value = '  +0.99%  '

$ws.Range("D51").Value = '''0.07161'
